# Insert a new data row at row 333 (shifting existing rows 333:360 down to 334:361)
# and populate it with a new "Acelga" price record, matching the author's weekly update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 333; this shifts rows 333-360 down to 334-361
# and copies formatting (incl. the date number format on column D) from the row above.
$ws.Rows.Item(333).Insert()

# Populate the newly inserted row 333 with the new record's values.
$ws.Range("A333").Value2 = 4
$ws.Range("B333").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C333").Value2 = "Los Lagos"
$ws.Range("D333").Value2 = 45223
$ws.Range("E333").Value2 = 10
$ws.Range("F333").Value2 = 100112009
$ws.Range("G333").Value2 = "Acelga"
$ws.Range("H333").Value2 = "Sin especificar"
$ws.Range("I333").Value2 = "Primera"
$ws.Range("J333").Value2 = 75
$ws.Range("K333").Value2 = 10000
$ws.Range("L333").Value2 = 10000
$ws.Range("M333").Value2 = 10000
$ws.Range("N333").Value2 = "`$/docena de atados (12 kilos)"
$ws.Range("O333").Value2 = "Región de La Araucanía"
$ws.Range("P333").Value2 = 833
$ws.Range("Q333").Value2 = 12
$ws.Range("R333").Value2 = "Hortaliza"
